$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells are stored as text so values like "1.00" or "3.20"
# keep their exact display form instead of being coerced to numbers.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "87.605.10"
$ws.Range("E2").Value = "  +0.31%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.178.52"
$ws.Range("E3").Value = "  -3.52%  "

# Row 4
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.46"
$ws.Range("E5").Value = "  -2.47%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "608.80"
$ws.Range("E6").Value = "  -2.85%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.388"
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.673"
$ws.Range("E8").Value = "  +5.41%  "

# Row 9
$ws.Range("E9").Value = "  -0.01%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.173.03"
$ws.Range("E10").Value = "  -3.65%  "

# Row 11
$ws.Range("E11").Value = "  -8.14%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.176"
$ws.Range("E12").Value = "  +0.80%  "

# Row 13
$ws.Range("E13").Value = "  -7.93%  "

# Row 14
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.758.36"
$ws.Range("E14").Value = "  -3.70%  "

# Row 15
$ws.Range("B15").Value = "Toncoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.28"
$ws.Range("E15").Value = "  +0.33%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.383.96"
$ws.Range("E16").Value = "  +0.62%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "32.21"
$ws.Range("E17").Value = "  -6.86%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.176.90"
$ws.Range("E18").Value = "  -3.07%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.20"
$ws.Range("E19").Value = "  +7.42%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.47"
$ws.Range("E20").Value = "  -5.52%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "413.23"
$ws.Range("E21").Value = "  -5.79%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.49"
$ws.Range("E22").Value = "  -7.93%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.08"
$ws.Range("E23").Value = "  -5.62%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.20"
$ws.Range("E24").Value = "  -0.52%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.18"
$ws.Range("E25").Value = "  -0.20%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.334.94"
$ws.Range("E26").Value = "  -3.98%  "

# Row 27
$ws.Range("E27").Value = "  +0.80%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "73.47"
$ws.Range("E28").Value = "  -4.61%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.19%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.162"
$ws.Range("E30").Value = "  -8.04%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.997"
$ws.Range("E31").Value = "  -0.41%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "549.70"
$ws.Range("E32").Value = "  -1.10%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.23"
$ws.Range("E33").Value = "  -7.94%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.32"
$ws.Range("E34").Value = "  -9.16%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.89"
$ws.Range("E35").Value = "  -0.35%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.86"
$ws.Range("E36").Value = "  -6.37%  "

# Row 37
$ws.Range("E37").Value = "  -6.00%  "

# Row 38
$ws.Range("E38").Value = "  -4.15%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "21.84"
$ws.Range("E39").Value = "  +0.35%  "

# Row 40
$ws.Range("E40").Value = "  -0.19%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.04"
$ws.Range("E41").Value = "  +1.99%  "

# Row 42
$ws.Range("E42").Value = "  -0.09%  "

# Row 43
$ws.Range("E43").Value = "  -5.00%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.370"
$ws.Range("E44").Value = "  -8.05%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "149.39"
$ws.Range("E45").Value = "  -2.67%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "173.68"
$ws.Range("E46").Value = "  -4.17%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "43.29"
$ws.Range("E47").Value = "  -3.75%  "

# Row 48
$ws.Range("E48").Value = "  +4.13%  "

# Row 49
$ws.Range("E49").Value = "  -8.86%  "

# Row 50
$ws.Range("E50").Value = "  -6.87%  "

# Row 51
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.74"
$ws.Range("E51").Value = "  -3.50%  "
